$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") all share the
# same layout: add a header label in A1 (styled like the year headers in
# B1:E1), drop the bold/bordered style from A2:A12 and fix a handful of
# accented labels.
# ---------------------------------------------------------------------------

$accentLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell A1, matching the style already used by B1:E1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Rows 2-12: drop the style and refresh the accented text.
    foreach ($r in 2..12) {
        $cell = $ws.Range("A" + $r)
        $cell.Value = $accentLabels[$r]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais (MtCO2eq)"): add header A1, fix accented labels,
# remove the style from A2:A3, and delete row 4 ("Teto").
# ---------------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").Style = "Normal"

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").Style = "Normal"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6 ("Custo Total (bilhões de R$)"): add header A1, change B1 label,
# fix accented labels in A2:A3 (dropping their style) and update values.
# ---------------------------------------------------------------------------

$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's text must stay a string ("2015"), not become a number: set it with a
# leading apostrophe (forces text), then restore the plain header style
# (copied from sheet 1's B1, which has the same base formatting without the
# quote-prefix flag that the apostrophe trick would otherwise leave behind).
$ws6.Range("B1").Value = "'2015"
$wb.Worksheets.Item(1).Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").Style = "Normal"
$ws6.Range("B2").Value = 573

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").Style = "Normal"
$ws6.Range("B3").Value = 99

Write-Output "edit applied"
